# The edit swaps the entire contents of data rows 2 and 3 (the two
# observation records got reordered). Most columns happen to hold the
# same value in both rows, so to avoid any unwanted side effects (e.g.
# Excel re-interpreting date/time text as real dates when a value is
# written back) we only touch the columns whose values actually differ
# between row 2 and row 3: A, B, D, E, F, G, H, Q, R, and AC.
#
# Column AC ("Publik kommentar") only has content in row 2 ("På sälg")
# before the edit; after the edit it only has content in row 3, so we
# clear it from row 2 and set it on row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$swapCols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $swapCols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")
    $v2 = $cell2.Value()
    $v3 = $cell3.Value()
    $cell2.Value = $v3
    $cell3.Value = $v2
}

# AC2 moves to AC3 ("På sälg" comment follows the record to its new row).
$ws.Range("AC2").ClearContents()
$ws.Range("AC3").Value = "På sälg"
